$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A82").Value = "3VAPLP"
$ws.Range("B82").Value = "Rodillo de separación Epson"
$ws.Range("C82").Value = "L1110 L1118 L1119 L1250 L3100 L3106 L3108 L3110 L3115 L3116 L3117 L3118 L3119 L3150 L3151 L3153 L3156 L3158 L3160 L3161 L3163 L3210 L3250 L4150 L4158 L4160 L4168 L4170 L4260 L5190 L6178"
$ws.Range("D82").Value = 0
$ws.Range("E82").Value = 100000
$ws.Range("F82").Value = 1
$ws.Range("G82").Value = 6
$ws.Range("H82").Formula = "=(E82-D82)*G82"
$ws.Range("I82").Formula = "=D82*F82"
$ws.Range("J82").Value = 0
